$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12 (shifts existing rows 12..58 down to 13..59)
$ws.Rows.Item(12).Insert()

# Populate the new row 12 with the new data record
$ws.Cells.Item(12, 1).Value = 4
$ws.Cells.Item(12, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(12, 3).Value = "Los Lagos"
$ws.Cells.Item(12, 4).Value = 44565
$ws.Cells.Item(12, 5).Value = 10
$ws.Cells.Item(12, 6).Value = "Fruta"
$ws.Cells.Item(12, 7).Value = 100103
$ws.Cells.Item(12, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(12, 9).Value = 100103001
$ws.Cells.Item(12, 10).Value = "Cereza"
$ws.Cells.Item(12, 11).Value = "Santina"
$ws.Cells.Item(12, 12).Value = "Primera"
$ws.Cells.Item(12, 13).Value = 800
$ws.Cells.Item(12, 14).Value = 7000
$ws.Cells.Item(12, 15).Value = 7500
$ws.Cells.Item(12, 16).Value = 7250
$ws.Cells.Item(12, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(12, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(12, 19).Value = 725
$ws.Cells.Item(12, 20).Value = 10
